$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "sdf"
$ws.Range("C3").Value = "sdagfafg"
$ws.Range("D3").Value = "asdfg"
$ws.Range("E3").Value = "asdfasdf"

$ws.Range("D5").Select()
